$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the first row's login labels (previously a name/city sample row)
$ws.Cells.Item(1, 1).Value2 = "USERNAME"
$ws.Cells.Item(1, 2).Value2 = "PASSWORD"

# Clear out the sample date values in column C, keeping their date formatting
$ws.Range("C1:C4").ClearContents()

# Match the author's final selection
$ws.Range("C2").Select()
